# Updating example upload files for v2.
#
# Insert a new "continuity_of_support" column into the Episodes sheet at
# column AB, pushing the existing "episode_tags" column from AB to AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Episodes")

# First copy the existing "episode_tags" header/value out to the new AC
# column (before AB gets overwritten), then write the new header+values
# into AB.
$ws.Range("AC1").Value = "episode_tags"
$ws.Range("AC2").Value = "tag3"

$ws.Range("AB1").Value = "continuity_of_support"
$ws.Range("AB2").Value = 2
$ws.Range("AB3").Value = 2
$ws.Range("AB4").Value = 1
$ws.Range("AB5").Value = 9

# Make Episodes the active/selected sheet with AC5 selected, matching the
# view state the author left the workbook in when saving.
$ws.Activate()
$ws.Range("AC5").Select()
